# Corrección a Diebold Mariano: actualiza matrices de p-valores y
# estadísticos DM en las hojas "P_valores" y "Estadisticos_DM".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Hoja "P_valores"
# ---------------------------------------------------------------
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.1299122391364895
$wsP.Range("D2").Value = 0.04891241990706074
$wsP.Range("E2").Value = 0.1410107461270433
$wsP.Range("F2").Value = 0.08089768487822124

$wsP.Range("B3").Value = 0.1299122391364895
$wsP.Range("D3").Value = 0.6787708002486448
$wsP.Range("E3").Value = 0.5068136697383347
$wsP.Range("F3").Value = 0.8651271143427492

$wsP.Range("B4").Value = 0.04891241990706074
$wsP.Range("C4").Value = 0.6787708002486448
$wsP.Range("E4").Value = 0.3285128584248436
$wsP.Range("F4").Value = 0.4988649131676739

$wsP.Range("B5").Value = 0.1410107461270433
$wsP.Range("C5").Value = 0.5068136697383347
$wsP.Range("D5").Value = 0.3285128584248436
$wsP.Range("F5").Value = 0.5230314518366534

$wsP.Range("B6").Value = 0.08089768487822124
$wsP.Range("C6").Value = 0.8651271143427492
$wsP.Range("D6").Value = 0.4988649131676739
$wsP.Range("E6").Value = 0.5230314518366534

# ---------------------------------------------------------------
# Hoja "Estadisticos_DM"
# ---------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -1.609066406760991
$wsE.Range("D2").Value = -2.156592795854996
$wsE.Range("E2").Value = -1.560270264753369
$wsE.Range("F2").Value = -1.881255762302034

$wsE.Range("B3").Value = 1.609066406760991
$wsE.Range("D3").Value = -0.4229298805370965
$wsE.Range("E3").Value = 0.681270547682354
$wsE.Range("F3").Value = 0.1730007348328235

$wsE.Range("B4").Value = 2.156592795854996
$wsE.Range("C4").Value = 0.4229298805370965
$wsE.Range("E4").Value = 1.012426389404873
$wsE.Range("F4").Value = 0.6942826192114511

$wsE.Range("B5").Value = 1.560270264753369
$wsE.Range("C5").Value = -0.681270547682354
$wsE.Range("D5").Value = -1.012426389404873
$wsE.Range("F5").Value = -0.6550852730162153

$wsE.Range("B6").Value = 1.881255762302034
$wsE.Range("C6").Value = -0.1730007348328235
$wsE.Range("D6").Value = -0.6942826192114511
$wsE.Range("E6").Value = 0.6550852730162153

$wb.Save()
